$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Text edits inside the two rich-text header cells (shared strings).
#    A8  = "Volume 31   Number  6"                         -> "...7"
#    C9  = "Report Covering the Week  2/5/2024  Through  2/11/2024"
#                                     -> 2/12/2024 ... 2/18/2024
# ---------------------------------------------------------------------------

$volRange = $ws.Range("A8")
$volText = $volRange.Value()
$volPos = $volText.IndexOf("6") + 1
$volRange.Characters($volPos, 1).Text = "7"

$weekRange = $ws.Range("C9")
$weekText = $weekRange.Value()
$startPos = $weekText.IndexOf("2/5/2024") + 1
$weekRange.Characters($startPos, 8).Text = "2/12/2024"

$weekText2 = $weekRange.Value()
$endPos = $weekText2.IndexOf("2/11/2024") + 1
$weekRange.Characters($endPos, 9).Text = "2/18/2024"

# ---------------------------------------------------------------------------
# 2) Style fix-ups: some cells change "kind" (blank placeholder text <->
#    real number) along with their value, which also changes which cell
#    style (number format) applies. Copying a same-style cell over first
#    guarantees the destination lands on the exact right style, and also
#    handles the text "0" case (plain Value assignment of a numeric-looking
#    string like "0" gets auto-coerced back into a number).
# ---------------------------------------------------------------------------

$textZeroStyleSrc = $ws.Range("C14")   # s=14, shared text "0"
$numStyle15Src    = $ws.Range("L15")   # s=15, numeric (percent-like)
$numStyle16Src    = $ws.Range("I22")   # s=16, numeric (plain count)

# Cells that need to become numeric with style 15 (were text placeholders)
foreach ($addr in @("N14","E15","H15","K15","E23","N28","N29")) {
    $numStyle15Src.Copy($ws.Range($addr))
}

# Cells that need to become numeric with style 16 (were text placeholders)
foreach ($addr in @("D15","G15","J15","D23")) {
    $numStyle16Src.Copy($ws.Range($addr))
}

# Cells that need to become the text "0" placeholder with style 14
foreach ($addr in @("C20","C22","F23","C28","C29")) {
    $textZeroStyleSrc.Copy($ws.Range($addr))
}

# ---------------------------------------------------------------------------
# 3) Write every updated numeric value (covers both the cells whose style
#    was just fixed above, and all the plain same-style value updates).
# ---------------------------------------------------------------------------

$values = @{
    "N14" = -100
    "D15" = 1
    "E15" = -100
    "G15" = 1
    "H15" = -100
    "J15" = 1
    "K15" = -100

    "C16" = 5
    "D16" = 6
    "E16" = -16.666666666666
    "F16" = 12
    "G16" = 15
    "H16" = -20
    "I16" = 22
    "J16" = 26
    "K16" = -15.384615384615
    "L16" = -4.347826086956
    "M16" = -4.347826086956
    "N16" = -81.355932203389

    "C17" = 3
    "D17" = 3
    "E17" = 0
    "F17" = 14
    "G17" = 14
    "I17" = 23
    "J17" = 29
    "K17" = -20.689655172413
    "L17" = -28.125
    "M17" = 64.285714285714
    "N17" = -25.806451612903

    "C18" = 5
    "D18" = 4
    "E18" = 25
    "F18" = 15
    "G18" = 19
    "H18" = -21.052631578947
    "I18" = 33
    "J18" = 33
    "K18" = 0
    "L18" = 3.125
    "M18" = -25
    "N18" = -87.912087912087

    "C19" = 10
    "D19" = 17
    "E19" = -41.176470588235
    "F19" = 71
    "G19" = 70
    "H19" = 1.428571428571
    "I19" = 129
    "J19" = 120
    "K19" = 7.5
    "L19" = 10.256410256410
    "M19" = -32.460732984293
    "N19" = -63.352272727272

    "D20" = 3
    "E20" = -100
    "F20" = 3
    "G20" = 4
    "H20" = -25
    "I20" = 4
    "J20" = 5
    "K20" = -20
    "L20" = -42.857142857142
    "M20" = -20
    "N20" = -97.790055248618

    "C21" = 23
    "D21" = 34
    "E21" = -32.352941176470
    "F21" = 115
    "G21" = 123
    "H21" = -6.504065040650
    "I21" = 211
    "J21" = 214
    "K21" = -1.401869158878
    "L21" = -1.401869158878
    "M21" = -24.100719424460
    "N21" = -77.997914494264

    "D22" = 4
    "E22" = -100
    "F22" = 6
    "G22" = 8
    "H22" = -25
    "J22" = 12
    "K22" = 16.666666666666
    "L22" = -17.647058823529

    "D23" = 1
    "E23" = -100
    "H23" = -100
    "J23" = 2
    "K23" = -50
    "M23" = -75

    "C24" = 61
    "D24" = 41
    "E24" = 48.780487804878
    "F24" = 222
    "H24" = 40.506329113924
    "I24" = 392
    "J24" = 273
    "K24" = 43.589743589743
    "L24" = 30.232558139534
    "M24" = 59.349593495935

    "C25" = 18
    "D25" = 13
    "E25" = 38.461538461538
    "F25" = 45
    "G25" = 32
    "H25" = 40.625
    "I25" = 68
    "J25" = 56
    "K25" = 21.428571428571
    "L25" = 28.301886792452
    "M25" = 41.666666666666

    "D26" = 2
    "G26" = 4
    "J26" = 4

    "F27" = 5
    "H27" = -16.666666666666
    "I27" = 13
    "J27" = 10
    "K27" = 30
    "L27" = 8.333333333333

    "N28" = 0
    "N29" = 0

    "I30" = 3
    "K30" = 200
    "L30" = 200
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
